$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

# Replace the placeholder contact data (Tom Peter / David Cris / Mukta Sharma)
# with the corrected names (Pooja Singh / Anita Singh / Ranvijay Singh).
$ws.Range("B2").Value = "Pooja"
$ws.Range("C2").Value = "Singh"
$ws.Range("B3").Value = "Anita"
$ws.Range("C3").Value = "Singh"
$ws.Range("B4").Value = "Ranvijay"
$ws.Range("C4").Value = "Singh"

# Move/restore the active selection to F5, matching the saved state.
$ws.Range("F5").Select()
